$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.750.91'
$ws.Range('E2').Value = '  +2.96%  '
$ws.Range('D3').Value = '3.133.81'
$ws.Range('E3').Value = '  +1.92%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.72'
$ws.Range('E5').Value = '  +1.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.44'
$ws.Range('E6').Value = '  +3.85%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.125.45'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.533'
$ws.Range('E9').Value = '  +1.23%  '
$ws.Range('E10').Value = '  +14.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.72'
$ws.Range('E11').Value = '  -0.62%  '
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('E13').Value = '  +5.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.45'
$ws.Range('E14').Value = '  +6.23%  '
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('D16').Value = '3.652.57'
$ws.Range('E16').Value = '  +1.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.18'
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('D18').Value = '63.654.31'
$ws.Range('E18').Value = '  +2.92%  '
$ws.Range('D19').Value = '3.131.24'
$ws.Range('E19').Value = '  +1.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '465.39'
$ws.Range('E20').Value = '  +4.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.39'
$ws.Range('E21').Value = '  +3.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.733'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('E23').Value = '  +1.68%  '
$ws.Range('E24').Value = '  -3.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.28'
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.98'
$ws.Range('E27').Value = '  +8.89%  '
$ws.Range('E28').Value = '  +1.80%  '
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.91'
$ws.Range('E31').Value = '  +2.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.17'
$ws.Range('E32').Value = '  +1.20%  '
$ws.Range('E33').Value = '  -2.06%  '
$ws.Range('D34').Value = '0.0₃0880'
$ws.Range('E34').Value = '  +11.17%  '
$ws.Range('E35').Value = '  +8.41%  '
$ws.Range('E36').Value = '  +1.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.42'
$ws.Range('E37').Value = '  +15.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.14'
$ws.Range('E38').Value = '  +1.54%  '
$ws.Range('E39').Value = '  +1.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '453.80'
$ws.Range('E40').Value = '  +7.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.75'
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('E42').Value = '  +0.96%  '
$ws.Range('D43').Value = '2.907.04'
$ws.Range('E43').Value = '  -1.26%  '
$ws.Range('E44').Value = '  +1.81%  '
$ws.Range('E45').Value = '  +1.88%  '
$ws.Range('E46').Value = '  +2.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '127.79'
$ws.Range('E47').Value = '  +2.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '35.94'
$ws.Range('E48').Value = '  +2.67%  '
$ws.Range('E50').Value = '  +0.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.74'
$ws.Range('E51').Value = '  +1.69%  '
